$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.412.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6308"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.911.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.010"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001088"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6810"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.185.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.209"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.465.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.525"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.361"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.300"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05584"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.037"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.857"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.588"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.249.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01810"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.431"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9030"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.157"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4017"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.693"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.983"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1123"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
